# Updates cryptos list: refresh Price/Volume(1h) figures for the
# existing rows and reorder rows 39-41 (Maker / VeChain / FraxShare)
# with their refreshed data, matching the Aug 17 2023 GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.471.44"
$ws.Range("E2").Value = "  -2.72%  "
$ws.Range("D3").Value = "1.774.89"
$ws.Range("E3").Value = "  -3.08%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.47"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5870"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2741"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.22"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("E10").Value = "  -4.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07531"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "1.785.54"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.749"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6067"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").Value = "2.015.93"
$ws.Range("E15").Value = "  -3.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "74.73"
$ws.Range("E16").Value = "  -4.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008614"
$ws.Range("E17").Value = "  -11.08%  "
$ws.Range("D18").Value = "28.430.95"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.356"
$ws.Range("E19").Value = "  -5.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "206.50"
$ws.Range("E21").Value = "  -6.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.36"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.729"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.40"
$ws.Range("E25").Value = "  -3.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.116"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1246"
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.25"
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.409"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06142"
$ws.Range("E30").Value = "  -4.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.411"
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.757"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.670"
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.043"
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6349"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.501"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.680"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01673"
$ws.Range("E39").Value = "  -5.04%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.298"
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.138.77"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8718"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.81"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "1.929.14"
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.54"
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000111"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.572"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.370"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05413"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4461"
$ws.Range("E51").Value = "  -2.11%  "
